# Add vaccination data sheet (Planilha4) to the workbook, mirroring the
# author's "start adding data about vaccination" commit.

$wb = $excel.ActiveWorkbook

# --- 1. Create the new worksheet as the last tab, named "Planilha4" -------
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Planilha4"

# --- 2. Borrow existing cell styles (so styles.xml gains no new records) --
# Style used by Planilha1!B1  -> centered, general number format (index 2)
$wb.Worksheets.Item("Planilha1").Range("B1").Copy() | Out-Null
$ws.Range("B1:E1").PasteSpecial(-4122) | Out-Null
$ws.Range("B2:E4").PasteSpecial(-4122) | Out-Null

# Style used by Planilha2!B1 -> centered + vcentered, general format (index 4)
$wb.Worksheets.Item("Planilha2").Range("B1").Copy() | Out-Null
$ws.Range("A1").PasteSpecial(-4122) | Out-Null

# Style used by Planilha2!A1 -> centered + vcentered, date format (index 3)
$wb.Worksheets.Item("Planilha2").Range("A1").Copy() | Out-Null
$ws.Range("A2:A264").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- 3. Header row ----------------------------------------------------
$ws.Range("A1").Value = "DATA"
$ws.Range("B1").Value = "DOSE_PRIMEIRA"
$ws.Range("C1").Value = "DOSE_SEGUNDA"
$ws.Range("D1").Value = "DOSE_UNICA"
$ws.Range("E1").Value = "TOTAL"

# --- 4. First data rows (with formulas) --------------------------------
$ws.Range("A2").Value = 44256
$ws.Range("B2").Formula = "=990+87+694"
$ws.Range("C2").Formula = "=536+80"
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 2387

$ws.Range("A3").Value = 44257
$ws.Range("B3").Formula = "=990+87+694"
$ws.Range("C3").Formula = "=536+80"
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 2387

$ws.Range("A4").Value = 44258
$ws.Range("B4").Formula = "=990+87+694"
$ws.Range("C4").Formula = "=541+80"
$ws.Range("D4").Value = 0
$ws.Range("E4").Formula = "=B4+C4+D4"

# --- 5. Remaining dates, column A only, rows 5-264 ----------------------
for ($r = 5; $r -le 264; $r++) {
    $ws.Cells.Item($r, 1).Value = 44256 + ($r - 2)
}

# --- 6. Column widths (best effort; exact bestFit pixel widths are a
#        desktop-Excel font-metric artifact this runtime cannot reproduce) -
$ws.Columns.Item(1).ColumnWidth = 16
$ws.Range("B1:C1").ColumnWidth = 15.42578125
$ws.Columns.Item(4).ColumnWidth = 12.42578125
$ws.Columns.Item(5).ColumnWidth = 15.140625
$ws.Columns.Item(6).ColumnWidth = 12.28515625

# --- 7. Selection / view state ------------------------------------------
$ws.Range("B5").Select()
